$wb = $excel.ActiveWorkbook

# --- Hoja1: update the conversion note text ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.94 = 10962.36 pesos`n✅ 10962.36 pesos = 2.93 = 938.74 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- tasas: update rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 339.799
$ws2.Range("N12").Value = 3737
$ws2.Range("O12").Value = 320.01
